$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The data refresh run swapped two clashing team-name rows ("Sava Strmec" /
# "NK Bistra") everywhere they occur, and additionally found that three
# fixture rows (8/9, 19/20, 111/112) had their Home/Away teams - and every
# other stat column alongside them - reported on the wrong source row, so
# those row pairs get their whole record (columns B..AB, except the Date in
# column D) swapped between the two rows.
# ---------------------------------------------------------------------------

# 1) Simple single-cell team-name corrections: every other remaining
#    occurrence of "Sava Strmec" / "NK Bistra" in Home/Away columns just
#    flips to the other name.
$flipCells = @(
    "F12", "E18", "F22", "F28", "F33", "E38", "F38", "E54", "F69",
    "E78", "E79", "F83", "E92", "E101", "E102", "E109", "E110",
    "E120", "F120", "E124", "F128", "E129", "E134"
)

foreach ($addr in $flipCells) {
    $cell = $ws.Range($addr)
    $cur = $cell.Value()
    if ($cur -eq "Sava Strmec") {
        $cell.Value = "NK Bistra"
    } elseif ($cur -eq "NK Bistra") {
        $cell.Value = "Sava Strmec"
    }
}

# 2) Full row-pair swaps: columns B through AB, skipping D (Date), swap the
#    content between the two rows of each pair.
$rowPairs = @(
    @(8, 9),
    @(19, 20),
    @(111, 112)
)

# Column B=2 .. AB=28 ; D=4 is skipped (Date stays put)
$cols = 2..28 | Where-Object { $_ -ne 4 }

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
